$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'27.331.21"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.Value = "'  +0.92%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'1.825.53"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.Value = "'  +0.00%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 5)
$cell.Value = "'  -0.13%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'314.76"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.Value = "'  +0.60%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 5)
$cell.Value = "'  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'0.4488"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.Value = "'  -1.75%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(8, 4)
$cell.Value = "'0.3780"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.Value = "'  +1.34%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 5)
$cell.Value = "'  +1.71%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'0.8868"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'20.96"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.Value = "'  +0.08%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.Value = "'1.823.89"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.Value = "'  -0.15%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 4)
$cell.Value = "'6.738"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.Value = "'  +0.58%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'5.458"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.Value = "'  +1.74%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'93.57"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.Value = "'  +0.59%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'0.07131"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.Value = "'  +0.46%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 5)
$cell.Value = "'  -0.12%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'0.000008805"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.Value = "'  -0.54%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.Value = "'15.15"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.Value = "'  +0.76%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'27.336.38"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.Value = "'  +0.81%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'5.397"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.Value = "'  +3.82%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'10.97"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.Value = "'  -0.11%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 2)
$cell.Value = "'Toncoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 3)
$cell.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'1.968"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.Value = "'  -1.68%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 2)
$cell.Value = "'Monero"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 3)
$cell.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.Value = "'151.57"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.Value = "'  -0.16%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 2)
$cell.Value = "'LidoDAOToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 3)
$cell.Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.Value = "'2.323"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.Value = "'  +4.48%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 2)
$cell.Value = "'EthereumClassic"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 3)
$cell.Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'18.70"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.Value = "'  +1.24%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 2)
$cell.Value = "'InternetComputer(DFINITY)"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 3)
$cell.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'5.392"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 5)
$cell.Value = "'  +2.13%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 2)
$cell.Value = "'BitcoinCash"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 3)
$cell.Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.Value = "'117.89"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.Value = "'  +0.41%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 2)
$cell.Value = "'Stellar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 3)
$cell.Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'0.08871"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.Value = "'  -0.22%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 2)
$cell.Value = "'ImmutableX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 3)
$cell.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'0.7924"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.Value = "'  +4.33%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 2)
$cell.Value = "'ARBITRUM"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 3)
$cell.Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'1.200"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.Value = "'  +0.18%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 2)
$cell.Value = "'Filecoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 3)
$cell.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'4.604"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.Value = "'  +2.89%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 2)
$cell.Value = "'HuobiToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 3)
$cell.Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'2.914"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.Value = "'  -1.96%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 2)
$cell.Value = "'Frax"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'1.000"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.Value = "'  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 2)
$cell.Value = "'TrustWalletToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.Value = "'1.111"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.Value = "'  +0.65%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 2)
$cell.Value = "'VeChain"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 3)
$cell.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'0.01981"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.Value = "'  +0.48%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 2)
$cell.Value = "'Hedera"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'0.05308"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.Value = "'  +0.36%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 2)
$cell.Value = "'FraxShare"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'7.312"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.Value = "'  +1.49%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 2)
$cell.Value = "'TheSandbox"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 3)
$cell.Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'0.5334"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.Value = "'  -0.72%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 2)
$cell.Value = "'MXToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.Value = "'2.870"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.Value = "'  -0.47%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 2)
$cell.Value = "'Algorand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'0.1715"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.Value = "'  -0.03%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 2)
$cell.Value = "'RenderToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 3)
$cell.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'2.319"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.Value = "'  +17.24%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 2)
$cell.Value = "'Aptos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 3)
$cell.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'8.658"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.Value = "'  +0.63%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 2)
$cell.Value = "'Decentraland"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 3)
$cell.Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'0.5078"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.Value = "'  -3.33%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 2)
$cell.Value = "'EnergySwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'10.66"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.Value = "'  -0.14%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 2)
$cell.Value = "'NEARProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.Value = "'1.699"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.Value = "'  +1.18%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 2)
$cell.Value = "'Quant"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 3)
$cell.Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.Value = "'105.28"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.Value = "'  -0.47%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 2)
$cell.Value = "'PaxDollar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 3)
$cell.Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'1.000"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.Value = "'  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 2)
$cell.Value = "'Cronos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'0.06410"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.Value = "'  -0.07%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 2)
$cell.Value = "'Aave"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'66.02"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.Value = "'  +3.92%  "
$cell.Style = "Normal"
